$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# "b.md" has been handed off again: update its status everywhere,
# record the new handoff package + timestamp, and flag that the
# handback file is stale (commit: "Generate Report for Handoff").
# ------------------------------------------------------------------

$newStatus       = "Ready for handoff"
$newHandoffTime  = "2016-08-17 14:37:15"
$errorDetail     = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8ef6b50223b31b38ebce930c174e8154c45d73f6/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/58370dc20bbc3e5d7a3abfb505bf2d48efd44f86/e2e/b.md."

# ---- Overview sheet: row 3 is b.md (A=File Name, E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date) ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Cells.Item(3, 5).Value = $newStatus       # E3 zh-cn
$wsOverview.Cells.Item(3, 6).Value = $newStatus       # F3 de-de
$wsOverview.Cells.Item(3, 7).Value = $newHandoffTime  # G3 Latest HO Xliff Generate Date

# ---- zh-cn sheet: row 3 is b.md ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Cells.Item(3, 3).Value  = $newStatus                                                        # C3 Status
$wsZhCn.Cells.Item(3, 7).Value  = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"             # G3 Latest Handoff File
$wsZhCn.Cells.Item(3, 8).Value  = "2016-08-17 14:37:03"                                             # H3 Latest Handoff Datetime
$wsZhCn.Cells.Item(3, 16).Value = $errorDetail                                                       # P3 Error Detail
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664                                           # col P width -> 40

# ---- de-de sheet: row 3 is b.md ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Cells.Item(3, 3).Value  = $newStatus                                                        # C3 Status
$wsDeDe.Cells.Item(3, 7).Value  = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"             # G3 Latest Handoff File
$wsDeDe.Cells.Item(3, 8).Value  = $newHandoffTime                                                    # H3 Latest Handoff Datetime
$wsDeDe.Cells.Item(3, 16).Value = $errorDetail                                                       # P3 Error Detail
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664                                           # col P width -> 40
